# The deck originally ships with two theme parts:
#   ppt/theme/theme1.xml  -> "Integral"      (used by the slide master / all slides)
#   ppt/theme/theme2.xml  -> "Office Theme"  (used by the notes master)
#
# The target edit swaps them: the deck's visible theme (theme1.xml, reached
# through the slide master) becomes the stock "Office Theme" colour set,
# while the notes-only theme becomes "Integral". The PowerPoint object model
# exposes a theme's 12 colour-scheme slots as ThemeColorScheme.Item(1..12),
# each a simple RGB value, so we recolour the master's theme to the Office
# palette (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink, in that order).

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme

# Cosmetic rename attempts (theme / colour-scheme display names) - harmless
# no-ops if the host doesn't expose them as settable.
$theme.Name = "Office Theme"
$theme.ThemeColorScheme.Name = "Office"

$colorScheme = $theme.ThemeColorScheme

# RGB() isn't available in this host, so pass the packed 0x00BBGGRR values
# PowerPoint's ColorFormat.RGB already uses.
$colorScheme.Item(1).RGB  = 0        # dk1      -> 000000
$colorScheme.Item(2).RGB  = 16777215 # lt1      -> FFFFFF
$colorScheme.Item(3).RGB  = 6968388  # dk2      -> 44546A
$colorScheme.Item(4).RGB  = 15132391 # lt2      -> E7E6E6
$colorScheme.Item(5).RGB  = 13998939 # accent1  -> 5B9BD5
$colorScheme.Item(6).RGB  = 3243501  # accent2  -> ED7D31
$colorScheme.Item(7).RGB  = 10855845 # accent3  -> A5A5A5
$colorScheme.Item(8).RGB  = 49407    # accent4  -> FFC000
$colorScheme.Item(9).RGB  = 12874308 # accent5  -> 4472C4
$colorScheme.Item(10).RGB = 4697456  # accent6  -> 70AD47
$colorScheme.Item(11).RGB = 12673797 # hlink    -> 0563C1
$colorScheme.Item(12).RGB = 7491477  # folHlink -> 954F72
